# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Estado de Cuenta" detail table (B16:G29) is rebuilt: for each
# worker the periods are re-listed in descending order (2403 -> 2309)
# and the two workers' blocks are grouped together (JUAN CAMILO PAJARO
# PEREZ first, then ALEJANDRO CORREA CUADRADO) instead of being
# interleaved period-by-period. The underlying (worker, period) ->
# (Valor Mora, Salario Basico) values are unchanged; only the row
# order/grouping changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# r, Tipo Doc, N Doc, Nombre, Periodo, Valor Mora, Salario Basico
$data = @(
  ,@(16, "CC", "1044933954", "JUAN CAMILO PAJARO PEREZ", "2403", 45600,  908526)
  ,@(17, "CC", "1044933954", "JUAN CAMILO PAJARO PEREZ", "2402", 72000,  908526)
  ,@(18, "CC", "1044933954", "JUAN CAMILO PAJARO PEREZ", "2401", 72000,  908526)
  ,@(19, "CC", "1044933954", "JUAN CAMILO PAJARO PEREZ", "2312", 72000,  908526)
  ,@(20, "CC", "1044933954", "JUAN CAMILO PAJARO PEREZ", "2311", 72000,  908526)
  ,@(21, "CC", "1044933954", "JUAN CAMILO PAJARO PEREZ", "2310", 72000,  908526)
  ,@(22, "CC", "1044933954", "JUAN CAMILO PAJARO PEREZ", "2309", 72000,  908526)
  ,@(23, "CC", "1051451471", "ALEJANDRO CORREA CUADRADO", "2403", 68400,  1400000)
  ,@(24, "CC", "1051451471", "ALEJANDRO CORREA CUADRADO", "2402", 108000, 1400000)
  ,@(25, "CC", "1051451471", "ALEJANDRO CORREA CUADRADO", "2401", 108000, 1400000)
  ,@(26, "CC", "1051451471", "ALEJANDRO CORREA CUADRADO", "2312", 108000, 1400000)
  ,@(27, "CC", "1051451471", "ALEJANDRO CORREA CUADRADO", "2311", 108000, 1400000)
  ,@(28, "CC", "1051451471", "ALEJANDRO CORREA CUADRADO", "2310", 108000, 1400000)
  ,@(29, "CC", "1051451471", "ALEJANDRO CORREA CUADRADO", "2309", 93600,  1400000)
)

foreach ($item in $data) {
  $r = $item[0]
  $ws.Cells.Item($r, 2).Value = $item[1]
  $ws.Cells.Item($r, 3).Value = $item[2]
  $ws.Cells.Item($r, 4).Value = $item[3]
  $ws.Cells.Item($r, 5).Value = $item[4]
  $ws.Cells.Item($r, 6).Value = $item[5]
  $ws.Cells.Item($r, 7).Value = $item[6]
}
